{"js": "// Update the two-digit multiplication expressions in the practice table.\n// Each \"before\" expression appears exactly once in the document, so a\n// straightforward search + replace per pair is unambiguous.\nconst replacements = [\n  [\"32\u00d781=\", \"15\u00d739=\"],\n  [\"43\u00d747=\", \"70\u00d727=\"],\n  [\"13\u00d711=\", \"25\u00d790=\"],\n  [\"52\u00d712=\", \"37\u00d760=\"],\n  [\"62\u00d785=\", \"44\u00d798=\"],\n  [\"50\u00d783=\", \"38\u00d711=\"],\n  [\"67\u00d758=\", \"46\u00d797=\"],\n  [\"65\u00d769=\", \"59\u00d724=\"],\n  [\"91\u00d733=\", \"64\u00d785=\"],\n  [\"80\u00d788=\", \"74\u00d743=\"],\n  [\"77\u00d739=\", \"76\u00d773=\"],\n  [\"25\u00d714=\", \"91\u00d732=\"],\n  [\"98\u00d742=\", \"70\u00d782=\"],\n  [\"76\u00d732=\", \"54\u00d776=\"],\n  [\"13\u00d766=\", \"12\u00d741=\"],\n  [\"55\u00d783=\", \"53\u00d769=\"],\n  [\"92\u00d787=\", \"79\u00d721=\"],\n  [\"98\u00d750=\", \"56\u00d794=\"],\n  [\"54\u00d763=\", \"88\u00d782=\"],\n  [\"13\u00d721=\", \"14\u00d785=\"],\n  [\"70\u00d791=\", \"89\u00d768=\"],\n  [\"44\u00d771=\", \"99\u00d752=\"],\n  [\"28\u00d716=\", \"66\u00d783=\"],\n  [\"46\u00d776=\", \"33\u00d719=\"],\n  [\"68\u00d730=\", \"13\u00d751=\"],\n];\n\nconst body = context.document.body;\n\nfor (const [oldText, newText] of replacements) {\n  const results = body.search(oldText, { matchCase: true, matchWholeWord: false });\n  results.load(\"items\");\n  await context.sync();\n\n  for (let i = 0; i < results.items.length; i++) {\n    results.items[i].insertText(newText, Word.InsertLocation.replace);\n  }\n  await context.sync();\n}\n", "ps1": "# Update the two-digit multiplication expressions in the practice table.\n# Each \"before\" expression appears exactly once in the document, so a\n# Find/Replace pass per pair is unambiguous.\n$d = $word.ActiveDocument\n\n$wdReplaceAll = 2\n$wdFindContinue = 1\n\n$replacements = @(\n    @(\"32\u00d781=\", \"15\u00d739=\"),\n    @(\"43\u00d747=\", \"70\u00d727=\"),\n    @(\"13\u00d711=\", \"25\u00d790=\"),\n    @(\"52\u00d712=\", \"37\u00d760=\"),\n    @(\"62\u00d785=\", \"44\u00d798=\"),\n    @(\"50\u00d783=\", \"38\u00d711=\"),\n    @(\"67\u00d758=\", \"46\u00d797=\"),\n    @(\"65\u00d769=\", \"59\u00d724=\"),\n    @(\"91\u00d733=\", \"64\u00d785=\"),\n    @(\"80\u00d788=\", \"74\u00d743=\"),\n    @(\"77\u00d739=\", \"76\u00d773=\"),\n    @(\"25\u00d714=\", \"91\u00d732=\"),\n    @(\"98\u00d742=\", \"70\u00d782=\"),\n    @(\"76\u00d732=\", \"54\u00d776=\"),\n    @(\"13\u00d766=\", \"12\u00d741=\"),\n    @(\"55\u00d783=\", \"53\u00d769=\"),\n    @(\"92\u00d787=\", \"79\u00d721=\"),\n    @(\"98\u00d750=\", \"56\u00d794=\"),\n    @(\"54\u00d763=\", \"88\u00d782=\"),\n    @(\"13\u00d721=\", \"14\u00d785=\"),\n    @(\"70\u00d791=\", \"89\u00d768=\"),\n    @(\"44\u00d771=\", \"99\u00d752=\"),\n    @(\"28\u00d716=\", \"66\u00d783=\"),\n    @(\"46\u00d776=\", \"33\u00d719=\"),\n    @(\"68\u00d730=\", \"13\u00d751=\")\n)\n\nforeach ($pair in $replacements) {\n    $oldText = $pair[0]\n    $newText = $pair[1]\n\n    $find = $d.Content.Find\n    $find.ClearFormatting()\n    $find.Replacement.ClearFormatting()\n    $find.Text = $oldText\n    $find.Replacement.Text = $newText\n    $find.Forward = $true\n    $find.Wrap = $wdFindContinue\n    $find.Execute($oldText, $false, $false, $false, $false, $false, $true, $wdFindContinue, $false, $newText, $wdReplaceAll)\n}\n"}
